$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: Dency Patel ---
$ws.Range("A2").Value = "Dency Patel"
$ws.Range("B2").Value = "dency.patel@accolitedigital.com"
$ws.Range("C2").Value = "'7629163826"
$ws.Range("D2").Value = "REACT"
$ws.Range("E2").Value = "Date:2023-04-26    FromTime:10:00:00     ToTime:17:00:00`nDate:2023-04-27    FromTime:10:00:00     ToTime:17:00:00"

# --- Row 3: Pratyush Singh ---
$ws.Range("A3").Value = "Pratyush Singh"
$ws.Range("B3").Value = "pratyush.singh@accolitedigital.com"
$ws.Range("C3").Value = "'8761826384"
$ws.Range("D3").Value = "SPRING"
$ws.Range("E3").Value = "Date:2023-04-28    FromTime:10:00:00     ToTime:17:00:00`nDate:2023-04-29    FromTime:09:00:00     ToTime:17:00:00"

# --- Remove old rows 4, 5, 6 (trainers no longer present) ---
$ws.Range("A4:E6").EntireRow.Delete()
